$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.258.91"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.446.60"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.20%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "581.92"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.27%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.95"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").Value = "2.441.49"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  -2.31%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.45"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").Value = "2.895.63"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "62.192.95"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "2.444.08"
$ws.Range("E18").Value = "  +0.65%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.79"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.26%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.12"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "327.46"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  -0.03%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "65.78"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.09%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.22"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.75%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "591.42"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -4.19%  "
$ws.Range("E28").Value = "  +1.86%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("E32").Value = "  -0.67%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.89"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.377"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("E39").Value = "  +4.53%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.41"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.55%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.78%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "43.29"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +1.87%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "142.49"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "0.0₆0265"
$ws.Range("E48").Value = "  +20.42%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.604"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.99%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0522"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "19.90"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
